$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the "QUERY" label cell) so the QUERY(...) formula in row 2
# shifts up into A1 -- matches the diff: row1/A1("QUERY") removed, row2
# becomes row1 keeping the same formula, dimension shrinks to A1, and the
# selection becomes the full-row selection left behind by a row delete.
$ws.Rows(1).Delete()

# Mirror the leftover full-row selection that Excel shows after deleting a
# row (selection collapses to the row that slid up into the deleted row's
# place): sqref="A1:XFD1".
[void]$ws.Rows(1).Select()
